$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 486; this shifts existing rows 486-593 down to 487-594
$ws.Rows.Item(486).Insert()

# Fill in the new row 486 with the new data (columns A-T)
$ws.Cells.Item(486, 1).Value = 6
$ws.Cells.Item(486, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(486, 3).Value = "Metropolitana"
$ws.Cells.Item(486, 4).Value = 44504
$ws.Cells.Item(486, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(486, 5).Value = 13
$ws.Cells.Item(486, 6).Value = "Fruta"
$ws.Cells.Item(486, 7).Value = 100101
$ws.Cells.Item(486, 8).Value = "Berries"
$ws.Cells.Item(486, 9).Value = 100101007
$ws.Cells.Item(486, 10).Value = "Kiwi"
$ws.Cells.Item(486, 11).Value = "Hayward"
$ws.Cells.Item(486, 12).Value = "Extra (doble especial)"
$ws.Cells.Item(486, 13).Value = 14
$ws.Cells.Item(486, 14).Value = 430000
$ws.Cells.Item(486, 15).Value = 450000
$ws.Cells.Item(486, 16).Value = 440000
$ws.Cells.Item(486, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(486, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(486, 19).Value = 978
$ws.Cells.Item(486, 20).Value = 450
